# Plantilla Visitas - actualizacion semanal de ingresos
# - limpia las patentes vehiculares ya usadas en las filas 7-9
# - corre la "Fecha de ingreso" una semana (19-dic-2025 -> 26-dic-2025)
# - ensancha un poco las columnas de RUT y Patente Vehicular
# - subraya la celda H15 (donde quedo el cursor) y deja la seleccion ahi

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Las patentes de las filas 7, 8 y 9 ya no corresponden: se borran
$ws.Range("G7").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("G9").ClearContents()

# La fecha de ingreso pasa de 19/12/2025 a 26/12/2025 para las filas con datos
$ws.Range("I4:I9").Value = 46017

# Las columnas D (RUT) y G (Patente Vehicular) quedan un poco mas anchas
$ws.Columns("D").ColumnWidth = 13.3
$ws.Columns("G").ColumnWidth = 9.1

# Se resalta con subrayado la celda H15 y se deja seleccionada
$ws.Range("H15").Font.Underline = 2
$ws.Range("H15").Select() | Out-Null
